$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update existing reason text (row 3 / A3) ---
$ws.Range("A3").Value = "FSSAI License Image Missing"

# --- Remove stray leading space from "Cropped Image" reason (row 26 / A26) ---
$ws.Range("A26").Value = "Cropped Image "

# --- Append three new Reason / Justification rows ---
$ws.Range("A58").Value = "Clear Display of Puzzle Pieces in Product Images"
$ws.Range("B58").Value = "High-quality images must clearly display the actual puzzle pieces included in the jigsaw puzzle pack. The visuals should illustrate the quantity, layout, and arrangement of the pieces to accurately represent how they are presented in the set. This will help customers better understand the contents, improve transparency, and minimize post-purchase confusion or returns."

$ws.Range("A59").Value = "Size Reference Missing"
$ws.Range("B59").Value = "The product listing is missing a required size reference. Please either provide a size guide or  include clear product dimensions (e.g., diameter in inches or centimeters), or explicitly state the size (e.g., Small, Medium, Large, or weight range) in the product title. This information is mandatory for catalog accuracy, searchability, and customer clarity."

$ws.Range("A60").Value = "Image Quality Issue"
$ws.Range("B60").Value = "The submitted product image of the back panel packaging is not clearly readable due to poor resolution, low contrast, or inadequate lighting. Key information such as ingredients, nutritional facts, allergen warnings, and usage instructions are difficult to discern, which compromises accuracy in catalog documentation and customer safety.All product images  must be high-resolution, sharp, and fully legible. Please resubmit a clear, well-lit image of the back panel packaging where all text and details are easily readable."

# --- Update sheet view (active selection) to match saved state ---
$ws.Range("A63").Select()
